$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(5, 7).Value = 1.91   # G5: 1.85 -> 1.91
$ws.Cells.Item(5, 9).Value = 3.8   # I5: 3.9 -> 3.8
$ws.Cells.Item(5, 10).Value = 2.5   # J5: 2.4 -> 2.5
$ws.Cells.Item(5, 11).Value = 2.3   # K5: 2.38 -> 2.3
$ws.Cells.Item(5, 15).Value = 1.22   # O5: 1.2 -> 1.22
$ws.Cells.Item(5, 16).Value = 4   # P5: 4.33 -> 4
$ws.Cells.Item(5, 17).Value = 1.73   # Q5: 1.7 -> 1.73
$ws.Cells.Item(5, 18).Value = 2.08   # R5: 2.1 -> 2.08
$ws.Cells.Item(5, 19).Value = 1.33   # S5: 1.3 -> 1.33
$ws.Cells.Item(5, 20).Value = 3.25   # T5: 3.4 -> 3.25
$ws.Cells.Item(5, 27).Value = 15   # AA5: 13 -> 15
$ws.Cells.Item(5, 29).Value = 13   # AC5: 15 -> 13
$ws.Cells.Item(5, 30).Value = 7   # AD5: 7.5 -> 7
$ws.Cells.Item(5, 34).Value = 13   # AH5: 15 -> 13
$ws.Cells.Item(5, 39).Value = 34   # AM5: 29 -> 34
$ws.Cells.Item(5, 41).Value = 10   # AO5: 9.5 -> 10
$ws.Cells.Item(5, 42).Value = 19   # AP5: 17 -> 19
$ws.Cells.Item(5, 43).Value = 34   # AQ5: 29 -> 34
$ws.Cells.Item(5, 44).Value = 51   # AR5: 41 -> 51
$ws.Cells.Item(5, 46).Value = 3.25   # AT5: 3.4 -> 3.25
$ws.Cells.Item(5, 50).Value = 5.5   # AX5: 6 -> 5.5
$ws.Cells.Item(5, 51).Value = 19   # AY5: 21 -> 19
$ws.Cells.Item(5, 53).Value = 51   # BA5: 67 -> 51
$ws.Cells.Item(6, 7).Value = 1.53   # G6: 1.5 -> 1.53
$ws.Cells.Item(6, 10).Value = 2.05   # J6: 2 -> 2.05
$ws.Cells.Item(6, 11).Value = 2.4   # K6: 2.5 -> 2.4
$ws.Cells.Item(6, 15).Value = 1.2   # O6: 1.18 -> 1.2
$ws.Cells.Item(6, 16).Value = 4.33   # P6: 4.5 -> 4.33
$ws.Cells.Item(6, 17).Value = 1.65   # Q6: 1.6 -> 1.65
$ws.Cells.Item(6, 18).Value = 2.2   # R6: 2.3 -> 2.2
$ws.Cells.Item(6, 19).Value = 1.3   # S6: 1.29 -> 1.3
$ws.Cells.Item(6, 20).Value = 3.4   # T6: 3.5 -> 3.4
$ws.Cells.Item(6, 23).Value = 8   # W6: 8.5 -> 8
$ws.Cells.Item(6, 28).Value = 23   # AB6: 21 -> 23
$ws.Cells.Item(6, 31).Value = 17   # AE6: 15 -> 17
$ws.Cells.Item(6, 32).Value = 51   # AF6: 41 -> 51
$ws.Cells.Item(6, 34).Value = 17   # AH6: 19 -> 17
$ws.Cells.Item(6, 46).Value = 3.4   # AT6: 3.5 -> 3.4
$ws.Cells.Item(6, 52).Value = 34   # AZ6: 29 -> 34
$ws.Cells.Item(7, 7).Value = 3.6   # G7: 3.7 -> 3.6
$ws.Cells.Item(7, 8).Value = 3.8   # H7: 3.9 -> 3.8
$ws.Cells.Item(7, 9).Value = 1.9   # I7: 1.85 -> 1.9
$ws.Cells.Item(7, 13).Value = 1.04   # M7: 1.03 -> 1.04
$ws.Cells.Item(7, 14).Value = 13   # N7: 15 -> 13
$ws.Cells.Item(7, 17).Value = 1.73   # Q7: 1.7 -> 1.73
$ws.Cells.Item(7, 18).Value = 2.08   # R7: 2.1 -> 2.08
$ws.Cells.Item(7, 23).Value = 12   # W7: 13 -> 12
$ws.Cells.Item(7, 29).Value = 13   # AC7: 15 -> 13
$ws.Cells.Item(7, 31).Value = 13   # AE7: 15 -> 13
$ws.Cells.Item(7, 34).Value = 8.5   # AH7: 9 -> 8.5
$ws.Cells.Item(7, 38).Value = 15   # AL7: 13 -> 15
$ws.Cells.Item(7, 39).Value = 23   # AM7: 21 -> 23
$ws.Cells.Item(7, 42).Value = 26   # AP7: 23 -> 26
$ws.Cells.Item(7, 48).Value = 51   # AV7: 41 -> 51
$ws.Cells.Item(7, 49).Value = 501   # AW7: 451 -> 501
$ws.Cells.Item(7, 51).Value = 10   # AY7: 9.5 -> 10
$ws.Cells.Item(7, 52).Value = 19   # AZ7: 17 -> 19
$ws.Cells.Item(7, 53).Value = 34   # BA7: 29 -> 34
$ws.Cells.Item(7, 54).Value = 51   # BB7: 41 -> 51
$ws.Cells.Item(10, 7).Value = 2.05   # G10: 2.1 -> 2.05
$ws.Cells.Item(10, 8).Value = 3.7   # H10: 3.6 -> 3.7
$ws.Cells.Item(10, 10).Value = 2.6   # J10: 2.63 -> 2.6
$ws.Cells.Item(10, 15).Value = 1.17   # O10: 1.18 -> 1.17
$ws.Cells.Item(10, 16).Value = 5   # P10: 4.5 -> 5
$ws.Cells.Item(10, 19).Value = 1.29   # S10: 1.3 -> 1.29
$ws.Cells.Item(10, 20).Value = 3.5   # T10: 3.4 -> 3.5
$ws.Cells.Item(10, 26).Value = 19   # Z10: 21 -> 19
$ws.Cells.Item(10, 30).Value = 7.5   # AD10: 7 -> 7.5
$ws.Cells.Item(10, 46).Value = 3.5   # AT10: 3.4 -> 3.5
$ws.Cells.Item(11, 7).Value = 2   # G11: 2.3 -> 2
$ws.Cells.Item(11, 8).Value = 3.9   # H11: 3.8 -> 3.9
$ws.Cells.Item(11, 9).Value = 3.25   # I11: 2.7 -> 3.25
$ws.Cells.Item(11, 10).Value = 2.5   # J11: 2.75 -> 2.5
$ws.Cells.Item(11, 12).Value = 3.4   # L11: 3 -> 3.4
$ws.Cells.Item(11, 13).Value = 1.01   # M11: 1.02 -> 1.01
$ws.Cells.Item(11, 14).Value = 26   # N11: 21 -> 26
$ws.Cells.Item(11, 21).Value = 1.36   # U11: 1.33 -> 1.36
$ws.Cells.Item(11, 22).Value = 3   # V11: 3.25 -> 3
$ws.Cells.Item(11, 23).Value = 15   # W11: 17 -> 15
$ws.Cells.Item(11, 24).Value = 15   # X11: 17 -> 15
$ws.Cells.Item(11, 25).Value = 9.5   # Y11: 10 -> 9.5
$ws.Cells.Item(11, 26).Value = 21   # Z11: 26 -> 21
$ws.Cells.Item(11, 27).Value = 13   # AA11: 15 -> 13
$ws.Cells.Item(11, 29).Value = 26   # AC11: 23 -> 26
$ws.Cells.Item(11, 30).Value = 9   # AD11: 8.5 -> 9
$ws.Cells.Item(11, 31).Value = 11   # AE11: 10 -> 11
$ws.Cells.Item(11, 32).Value = 26   # AF11: 23 -> 26
$ws.Cells.Item(11, 35).Value = 23   # AI11: 21 -> 23
$ws.Cells.Item(11, 36).Value = 12   # AJ11: 11 -> 12
$ws.Cells.Item(11, 37).Value = 41   # AK11: 29 -> 41
$ws.Cells.Item(11, 38).Value = 21   # AL11: 19 -> 21
$ws.Cells.Item(11, 39).Value = 21   # AM11: 19 -> 21
$ws.Cells.Item(11, 40).Value = 4.75   # AN11: 5 -> 4.75
$ws.Cells.Item(11, 41).Value = 10   # AO11: 12 -> 10
$ws.Cells.Item(11, 43).Value = 29   # AQ11: 34 -> 29
$ws.Cells.Item(11, 50).Value = 6   # AX11: 5.5 -> 6
$ws.Cells.Item(11, 51).Value = 15   # AY11: 13 -> 15
$ws.Cells.Item(11, 54).Value = 51   # BB11: 41 -> 51
$ws.Cells.Item(19, 7).Value = 4.75   # G19: 4.5 -> 4.75
$ws.Cells.Item(19, 9).Value = 1.7   # I19: 1.73 -> 1.7
$ws.Cells.Item(19, 10).Value = 5   # J19: 4.75 -> 5
$ws.Cells.Item(19, 17).Value = 1.83   # Q19: 1.85 -> 1.83
$ws.Cells.Item(19, 18).Value = 2.03   # R19: 2 -> 2.03
$ws.Cells.Item(19, 27).Value = 41   # AA19: 34 -> 41
$ws.Cells.Item(19, 35).Value = 8   # AI19: 8.5 -> 8
$ws.Cells.Item(24, 7).Value = 3.2   # G24: 2.9 -> 3.2
$ws.Cells.Item(24, 9).Value = 2.2   # I24: 2.38 -> 2.2
$ws.Cells.Item(24, 10).Value = 4   # J24: 3.6 -> 4
$ws.Cells.Item(24, 12).Value = 3   # L24: 3.1 -> 3
$ws.Cells.Item(24, 23).Value = 9   # W24: 8.5 -> 9
$ws.Cells.Item(24, 25).Value = 12   # Y24: 11 -> 12
$ws.Cells.Item(24, 26).Value = 34   # Z24: 29 -> 34
$ws.Cells.Item(24, 27).Value = 29   # AA24: 26 -> 29
$ws.Cells.Item(24, 28).Value = 41   # AB24: 34 -> 41
$ws.Cells.Item(24, 33).Value = 351   # AG24: 301 -> 351
$ws.Cells.Item(24, 34).Value = 7   # AH24: 7.5 -> 7
$ws.Cells.Item(24, 35).Value = 10   # AI24: 11 -> 10
$ws.Cells.Item(24, 36).Value = 9.5   # AJ24: 10 -> 9.5
$ws.Cells.Item(24, 37).Value = 21   # AK24: 23 -> 21
$ws.Cells.Item(24, 38).Value = 19   # AL24: 21 -> 19
$ws.Cells.Item(24, 40).Value = 5   # AN24: 4.75 -> 5
$ws.Cells.Item(24, 41).Value = 19   # AO24: 17 -> 19
$ws.Cells.Item(24, 43).Value = 67   # AQ24: 51 -> 67
$ws.Cells.Item(24, 45).Value = 251   # AS24: 201 -> 251
$ws.Cells.Item(24, 47).Value = 8.5   # AU24: 8 -> 8.5
$ws.Cells.Item(24, 52).Value = 23   # AZ24: 26 -> 23
$ws.Cells.Item(24, 53).Value = 41   # BA24: 51 -> 41
$ws.Cells.Item(25, 7).Value = 3.4   # G25: 3 -> 3.4
$ws.Cells.Item(25, 8).Value = 3.3   # H25: 3.2 -> 3.3
$ws.Cells.Item(25, 9).Value = 2.05   # I25: 2.3 -> 2.05
$ws.Cells.Item(25, 10).Value = 3.75   # J25: 3.4 -> 3.75
$ws.Cells.Item(25, 12).Value = 2.63   # L25: 2.88 -> 2.63
$ws.Cells.Item(25, 21).Value = 1.57   # U25: 1.53 -> 1.57
$ws.Cells.Item(25, 22).Value = 2.25   # V25: 2.38 -> 2.25
$ws.Cells.Item(25, 23).Value = 13   # W25: 12 -> 13
$ws.Cells.Item(25, 24).Value = 19   # X25: 17 -> 19
$ws.Cells.Item(25, 25).Value = 12   # Y25: 11 -> 12
$ws.Cells.Item(25, 27).Value = 23   # AA25: 21 -> 23
$ws.Cells.Item(25, 28).Value = 29   # AB25: 26 -> 29
$ws.Cells.Item(25, 31).Value = 12   # AE25: 11 -> 12
$ws.Cells.Item(25, 32).Value = 41   # AF25: 34 -> 41
$ws.Cells.Item(25, 34).Value = 9.5   # AH25: 10 -> 9.5
$ws.Cells.Item(25, 35).Value = 11   # AI25: 13 -> 11
$ws.Cells.Item(25, 36).Value = 9   # AJ25: 9.5 -> 9
$ws.Cells.Item(25, 37).Value = 19   # AK25: 23 -> 19
$ws.Cells.Item(25, 38).Value = 15   # AL25: 17 -> 15
$ws.Cells.Item(25, 39).Value = 21   # AM25: 23 -> 21
$ws.Cells.Item(25, 40).Value = 5.5   # AN25: 5 -> 5.5
$ws.Cells.Item(25, 41).Value = 17   # AO25: 15 -> 17
$ws.Cells.Item(25, 42).Value = 23   # AP25: 21 -> 23
$ws.Cells.Item(25, 44).Value = 67   # AR25: 51 -> 67
$ws.Cells.Item(25, 47).Value = 7.5   # AU25: 7 -> 7.5
$ws.Cells.Item(25, 50).Value = 4.33   # AX25: 4.5 -> 4.33
$ws.Cells.Item(25, 51).Value = 11   # AY25: 12 -> 11
$ws.Cells.Item(25, 53).Value = 34   # BA25: 41 -> 34
$ws.Cells.Item(25, 55).Value = 101   # BC25: 126 -> 101
